$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B: existing B..H (姓名.. 邮箱) shift right to C..I,
# and the new column becomes the "company name" column.
$ws.Columns("B:B").Insert()

# New column B: header + sample value (company name used to match DB entry)
$ws.Range("B1").Value = "公司名称(与录入数据库名称一致)"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B2").Value = "阿里"

# Column widths (character units); engine pads with a constant +5/7 offset when
# serializing to the stored <col width> value, so back that out here.
$offset = 5/7
$ws.Columns("B:B").ColumnWidth = 28 - $offset
$ws.Columns("D:G").ColumnWidth = 22.875 - $offset
$ws.Columns("H:H").ColumnWidth = 19.375 - $offset
$ws.Columns("I:I").ColumnWidth = 17.625 - $offset

# The hyperlink on the "邮箱" (email) cell moved from H2 to I2 with the column
# insert; re-anchor it (Insert() does not move the hyperlink's own range ref).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 9), "mailto:aa@qq.com")
# Re-assert the original text-format style (Add() reformats the cell onto a
# fresh "hyperlink" style variant); this snaps it back onto the existing one.
$ws.Range("I2").NumberFormat = "@"

# Restore the current selection to match the saved workbook state.
$ws.Range("F13").Select()
